$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '69.638.76'
Set-TextCell 2 5 '  +3.16%  '
Set-TextCell 3 4 '3.373.35'
Set-TextCell 3 5 '  +4.52%  '
Set-TextCell 4 5 '  +0.02%  '
Set-TextCell 5 4 '192.72'
Set-TextCell 5 5 '  +5.78%  '
Set-TextCell 6 4 '593.54'
Set-TextCell 6 5 '  +2.37%  '
Set-TextCell 7 5 '  +0.09%  '
Set-TextCell 8 4 '0.607'
Set-TextCell 8 5 '  +1.02%  '
Set-TextCell 9 5 '  +3.41%  '
Set-TextCell 10 4 '6.75'
Set-TextCell 10 5 '  +3.22%  '
Set-TextCell 11 5 '  +2.80%  '
Set-TextCell 12 4 '3.966.54'
Set-TextCell 12 5 '  +4.73%  '
Set-TextCell 14 4 '28.75'
Set-TextCell 14 5 '  +3.78%  '
Set-TextCell 15 4 '69.671.59'
Set-TextCell 15 5 '  +3.13%  '
Set-TextCell 16 4 '0.0000172'
Set-TextCell 16 5 '  +2.37%  '
Set-TextCell 17 4 '3.376.26'
Set-TextCell 17 5 '  +5.06%  '
Set-TextCell 18 4 '451.10'
Set-TextCell 18 5 '  +14.44%  '
Set-TextCell 19 4 '5.85'
Set-TextCell 19 5 '  +1.95%  '
Set-TextCell 20 4 '13.79'
Set-TextCell 20 5 '  +2.65%  '
Set-TextCell 21 4 '7.82'
Set-TextCell 21 5 '  +3.77%  '
Set-TextCell 22 4 '73.53'
Set-TextCell 22 5 '  +3.77%  '
Set-TextCell 23 5 '  +0.09%  '
Set-TextCell 24 4 '3.514.75'
Set-TextCell 24 5 '  +4.37%  '
Set-TextCell 25 2 'Polygon'
Set-TextCell 25 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 25 4 '0.518'
Set-TextCell 25 5 '  +1.02%  '
Set-TextCell 26 2 'PEPE'
Set-TextCell 26 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 26 4 '0.0000122'
Set-TextCell 26 5 '  +3.89%  '
Set-TextCell 27 4 '0.193'
Set-TextCell 27 5 '  +4.08%  '
Set-TextCell 28 4 '9.60'
Set-TextCell 28 5 '  +0.49%  '
Set-TextCell 29 4 '0.994'
Set-TextCell 29 5 '  -0.43%  '
Set-TextCell 30 5 '  +2.77%  '
Set-TextCell 31 4 '23.26'
Set-TextCell 31 5 '  +2.77%  '
Set-TextCell 32 4 '5.63'
Set-TextCell 32 5 '  +0.90%  '
Set-TextCell 33 5 '  +3.91%  '
Set-TextCell 34 4 '7.05'
Set-TextCell 34 5 '  +0.80%  '
Set-TextCell 35 4 '0.998'
Set-TextCell 35 5 '  -0.05%  '
Set-TextCell 36 5 '  +4.03%  '
Set-TextCell 37 4 '164.85'
Set-TextCell 37 5 '  +1.88%  '
Set-TextCell 38 4 '1.93'
Set-TextCell 38 5 '  +2.78%  '
Set-TextCell 39 4 '27.20'
Set-TextCell 39 5 '  +3.42%  '
Set-TextCell 40 4 '0.823'
Set-TextCell 40 5 '  +2.35%  '
Set-TextCell 41 4 '4.60'
Set-TextCell 41 5 '  +0.86%  '
Set-TextCell 42 5 '  +0.49%  '
Set-TextCell 43 4 '2.746.77'
Set-TextCell 43 5 '  +5.35%  '
Set-TextCell 44 4 '2.54'
Set-TextCell 44 5 '  +3.50%  '
Set-TextCell 45 4 '25.60'
Set-TextCell 45 5 '  +4.02%  '
Set-TextCell 46 4 '0.0690'
Set-TextCell 46 5 '  +1.35%  '
Set-TextCell 47 4 '344.46'
Set-TextCell 47 5 '  +2.95%  '
Set-TextCell 48 4 '40.82'
Set-TextCell 48 5 '  +0.61%  '
Set-TextCell 49 5 '  +3.55%  '
Set-TextCell 50 4 '32.97'
Set-TextCell 50 5 '  +7.50%  '
Set-TextCell 51 5 '  +4.79%  '
